$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "(1.74)"
$ws.Range("B6").Value = "(0.64)"

$ws.Range("C4").Value = "(0.15)"
$ws.Range("C6").Value = "(0.48)"

$ws.Range("D4").Value = "(0.08)"
$ws.Range("D6").Value = "(0.3)"

$ws.Range("E4").Value = "(1.89)"
$ws.Range("E6").Value = "(1.28)"

$ws.Range("F4").Value = "(0.51)"
$ws.Range("F6").Value = "(1.47)"

$ws.Range("G4").Value = "(0.87)"
$ws.Range("G6").Value = "(1.06)"

$ws.Range("H4").Value = "(0.97)"
$ws.Range("H6").Value = "(1.31)"

$ws.Range("I4").Value = "(1.3)"
$ws.Range("I6").Value = "(1.26)"

$ws.Range("J4").Value = "(3.34)"
$ws.Range("J6").Value = "(1.74)"
